# Swap the species-record data between row 2 and row 3 on the active sheet.
# Columns involved: A (Id), B (Taxonsorteringsordning), E (TaxonId),
# F (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2

    $cellRow2.Value2 = $val3
    $cellRow3.Value2 = $val2
}
